$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values keyed by row number, regenerated from Strike# -> K conversion
$kValues = @{
    2 = 1
    3 = 0
    4 = 2
    5 = 1
    6 = 1
    7 = 2
    8 = 2
    9 = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 3
    14 = 3
    15 = 0
    16 = 2
    17 = 1
    18 = 2
    19 = 1
    20 = 2
    21 = 0
    22 = 2
    23 = 0
    24 = 3
    25 = 0
    26 = 3
    27 = 1
    28 = 2
    29 = 1
    30 = 1
    31 = 2
    32 = 2
    33 = 0
    34 = 1
    35 = 0
    36 = 0
    37 = 2
    38 = 1
    39 = 0
    40 = 1
    41 = 3
    42 = 1
    43 = 1
    44 = 0
    45 = 1
    46 = 3
    47 = 1
    48 = 2
    50 = 2
    51 = 1
    52 = 2
    53 = 1
    54 = 1
    55 = 2
    56 = 1
    57 = 1
    58 = 1
    59 = 2
    60 = 0
    61 = 1
    62 = 1
    64 = 1
    65 = 1
    67 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
